# Renamed Automator Page Object
#
# The "Objects" sheet holds an elementName/elementType field table used by
# an automation framework. This commit swaps the old person-profile style
# fields (firstName, secondName, age, address, education, experience,
# salary) for a new contact-style field set (contactName, code,
# mobileNumber, emailID, role, isTheUserAuthorizedSignatory) - all typed
# as "string" - and the table shrinks from 8 data rows to 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header: elementName / elementType) is unchanged.

# New field rows.
$ws.Range("A2").Value = "contactName"
$ws.Range("B2").Value = "string"

$ws.Range("A3").Value = "code"
$ws.Range("B3").Value = "string"

$ws.Range("A4").Value = "mobileNumber"
$ws.Range("B4").Value = "string"

$ws.Range("A5").Value = "emailID"
$ws.Range("B5").Value = "string"

$ws.Range("A6").Value = "role"
$ws.Range("B6").Value = "string"

$ws.Range("A7").Value = "isTheUserAuthorizedSignatory"
$ws.Range("B7").Value = "string"

# The table now only spans down to row 7 - remove the old row 8 entirely
# (it previously held "salary" / "float").
$ws.Rows(8).Delete()

# Deleting row 8 shrinks the pre-existing column-B list validation range
# (B2:B22 -> B2:B21) as a side effect, even though that validation itself
# is untouched by this commit. Put it back exactly as it was.
$validated = $ws.Range("B2:B22")
$validated.Validation.Delete()
$validated.Validation.Add(3, 1, 3, '"string,int,object,list_object,float,"', 0)
$validated.Validation.ShowInput = $false

# Update the selected cell to match the author's final cursor position.
$ws.Range("C5").Select()
